$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with freshly scraped values.
# D-column prices are forced to text (NumberFormat "@") so values like
# "617.60", "1.00" or "0.999" keep their exact displayed digits instead of
# being normalised by Excel general number parsing; the style is then reset
# back to Normal so no stray formatting is left behind.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "76.147.68"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.024.56"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("E4").Value = "  +0.02%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "197.63"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "617.60"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +3.94%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.30%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.207"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +6.97%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "3.022.45"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +3.86%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.440"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("E12").Value = "  -0.26%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "5.21"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +6.59%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "3.576.63"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "28.82"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +2.44%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "76.027.76"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.29%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.0000192"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +2.83%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.014.40"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +3.56%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.42"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "8.90"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +3.27%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "378.63"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +2.50%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "2.38"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +5.67%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.38"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.99%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "3.177.94"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +4.34%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "72.28"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +0.05%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "4.32"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.98%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.79"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  +2.22%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.28%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "8.25"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.39"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.21%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "491.47"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  +5.59%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +12.20%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "20.51"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.36%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "162.17"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.41%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "20.04"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +1.66%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "190.12"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +6.36%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.376"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -2.79%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -5.41%  "
$ws.Range("E43").Value = "  +0.03%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.11"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +4.84%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.778"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +19.65%  "
$ws.Range("E46").Value = "  +6.23%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "41.02"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  +5.70%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.591"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("E51").Value = "  +0.18%  "
